# Auto-generated edit script: adds Top3-artists-per-horoscope data to Sheet1 (G:Z), matching
# the "Personalities data with top3 artists information" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 'Albers, Josef'
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 'American'
$ws.Range("J2").Value = 1888
$ws.Range("K2").Value = 'Bottrop (M�_nster district, North Rhine-Westphalia, Germany)'
$ws.Range("L2").Value = 'Portfolio I, Folder 23'
$ws.Range("N2").Value = 'Lee, Edward B.'
$ws.Range("O2").Value = 72
$ws.Range("P2").Value = 'American'
$ws.Range("Q2").Value = 1876
$ws.Range("S2").Value = 'Sunrise 6:35 AM, Greenville, PA'
$ws.Range("U2").Value = 'Cole, Timothy|Century Company'
$ws.Range("V2").Value = 63
$ws.Range("W2").Value = 'American|American'
$ws.Range("X2").Value = 1852
$ws.Range("Y2").Value = 'London (Greater London, England, United Kingdom)|'
$ws.Range("Z2").Value = 'Madonna and Child, by Giovanni Bellini'

# Row 3
$ws.Range("G3").Value = 'Hancock, John'
$ws.Range("H3").Value = 61
$ws.Range("I3").Value = 'English'
$ws.Range("J3").Value = 1757
$ws.Range("K3").Value = 'England'
$ws.Range("L3").Value = 'Untitled (Buttercup, JH 4)'
$ws.Range("M3").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=111136&size=Medium'
$ws.Range("N3").Value = 'Munhall, Walter'
$ws.Range("O3").Value = 58
$ws.Range("P3").Value = 'American'
$ws.Range("Q3").Value = 1901
$ws.Range("S3").Value = 'Construction workers riveting a truss'
$ws.Range("T3").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=38731&size=Medium'
$ws.Range("U3").Value = 'Douden, Herbert C.'
$ws.Range("V3").Value = 42
$ws.Range("W3").Value = 'American'
$ws.Range("X3").Value = 1901
$ws.Range("Z3").Value = 'Walls; [interior elevations]'

# Row 4
$ws.Range("G4").Value = 'Swank, Luke'
$ws.Range("H4").Value = 338
$ws.Range("I4").Value = 'American'
$ws.Range("J4").Value = 1890
$ws.Range("K4").Value = 'Johnstown, PA'
$ws.Range("L4").Value = '(Fair)'
$ws.Range("N4").Value = 'Goya, Francisco de'
$ws.Range("O4").Value = 105
$ws.Range("P4").Value = 'Spanish'
$ws.Range("Q4").Value = 1746
$ws.Range("R4").Value = 'Fuendetodos (Zaragoza province, Aragon, Spain)'
$ws.Range("S4").Value = 'Aquellos polbos. (Those Specks of Dust.)'
$ws.Range("T4").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=94136&size=Medium'
$ws.Range("U4").Value = 'Jorn, Asger'
$ws.Range("V4").Value = 32
$ws.Range("W4").Value = 'Danish'
$ws.Range("X4").Value = 1914
$ws.Range("Y4").Value = '�rhus county (Denmark)'
$ws.Range("Z4").Value = 'Semantic Virility'

# Row 5
$ws.Range("G5").Value = 'Mills, Frederick P.'
$ws.Range("H5").Value = 224
$ws.Range("I5").Value = 'American'
$ws.Range("J5").Value = 1879
$ws.Range("L5").Value = 'Ceiling; [ornament drawings]'
$ws.Range("N5").Value = 'Hare, Clyde'
$ws.Range("O5").Value = 142
$ws.Range("P5").Value = 'American'
$ws.Range("Q5").Value = 1927
$ws.Range("R5").Value = 'Bloomington, Indiana'
$ws.Range("S5").Value = 'Automobiles on Liberty Bridge'
$ws.Range("T5").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=9880&size=Medium'
$ws.Range("U5").Value = 'Alechinsky, Pierre'
$ws.Range("V5").Value = 101
$ws.Range("W5").Value = 'Belgian'
$ws.Range("X5").Value = 1927
$ws.Range("Y5").Value = 'Brussels, Belgium'
$ws.Range("Z5").Value = 'Lino-Litho'

# Row 6
$ws.Range("G6").Value = 'Spruance, Benton M.'
$ws.Range("H6").Value = 128
$ws.Range("I6").Value = 'American'
$ws.Range("J6").Value = 1904
$ws.Range("K6").Value = 'Philadelphia, PA'
$ws.Range("L6").Value = 'Remainders'
$ws.Range("M6").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=126459&size=Medium'
$ws.Range("N6").Value = 'Bochner, Mel'
$ws.Range("O6").Value = 70
$ws.Range("P6").Value = 'American'
$ws.Range("Q6").Value = 1940
$ws.Range("R6").Value = 'Pittsburgh, Pennsylvania'
$ws.Range("S6").Value = 'Design for Kraus Campo, Carnegie Mellon University: Study for"You Can Call it That if You Like"'
$ws.Range("T6").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=133238&size=Medium'
$ws.Range("U6").Value = 'Callot, Jacques'
$ws.Range("V6").Value = 63
$ws.Range("W6").Value = 'French'
$ws.Range("X6").Value = 1592
$ws.Range("Y6").Value = 'Nancy (Meurthe-et-Moselle, Lorraine, France)'
$ws.Range("Z6").Value = 'Drill with the Musket'

# Row 7
$ws.Range("G7").Value = 'Kauffman, William'
$ws.Range("H7").Value = 90
$ws.Range("I7").Value = 'American'
$ws.Range("J7").Value = 1857
$ws.Range("L7").Value = 'Courthouse; Westmoreland County Courthouse, Greensburg, PA; [detail drawings, plan]'
$ws.Range("N7").Value = 'Hood, Samuel S.'
$ws.Range("O7").Value = 85
$ws.Range("P7").Value = 'American'
$ws.Range("Q7").Value = 1917
$ws.Range("S7").Value = 'Leonard Lieb (Face)'
$ws.Range("T7").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=111697&size=Medium'
$ws.Range("U7").Value = 'Lefebre, John'
$ws.Range("V7").Value = 77
$ws.Range("W7").Value = 'American'
$ws.Range("X7").Value = 1905
$ws.Range("Y7").Value = 'Berlin, Germany'
$ws.Range("Z7").Value = 'Asger Jorn (Albisola, 1955)'

# Row 8
$ws.Range("G8").Value = 'Smith, W. Eugene'
$ws.Range("H8").Value = 571
$ws.Range("I8").Value = 'American'
$ws.Range("J8").Value = 1918
$ws.Range("K8").Value = 'Wichita, Kansas'
$ws.Range("L8").Value = 'City Council Chamber, City Council Building'
$ws.Range("M8").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=28980&size=Medium'
$ws.Range("N8").Value = 'Whistler, James McNeill'
$ws.Range("O8").Value = 125
$ws.Range("P8").Value = 'American'
$ws.Range("Q8").Value = 1834
$ws.Range("R8").Value = 'Lowell, Massachusetts'
$ws.Range("S8").Value = 'The Little Nude Model, Reading'
$ws.Range("T8").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=56834&size=Medium'
$ws.Range("U8").Value = 'Rijn, Rembrandt H. Van'
$ws.Range("V8").Value = 64
$ws.Range("W8").Value = 'Dutch'
$ws.Range("X8").Value = 1606
$ws.Range("Y8").Value = 'Netherlands, Leyden'
$ws.Range("Z8").Value = 'Saints Peter and John Healing the Cripple at the Gate of the Temple'

# Row 9
$ws.Range("G9").Value = 'Ruzicka, Rudolph'
$ws.Range("H9").Value = 95
$ws.Range("I9").Value = 'American'
$ws.Range("J9").Value = 1883
$ws.Range("K9").Value = 'Bohemia (Czech Republic)'
$ws.Range("L9").Value = 'The Washington Monument in the Public Garden, Boston'
$ws.Range("M9").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=100978&size=Medium|http://www.cmoa.org/CollectionImage.aspx?irn=100979&size=Medium'
$ws.Range("N9").Value = 'Rouault, Georges'
$ws.Range("O9").Value = 88
$ws.Range("P9").Value = 'French'
$ws.Range("Q9").Value = 1871
$ws.Range("R9").Value = 'Paris, France'
$ws.Range("S9").Value = 'Woman with Necklace'
$ws.Range("T9").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=124800&size=Medium'
$ws.Range("U9").Value = 'Hassam, Childe'
$ws.Range("V9").Value = 83
$ws.Range("W9").Value = 'American'
$ws.Range("X9").Value = 1859
$ws.Range("Y9").Value = 'Dorchester (Suffolk county, Massachusetts, United States)'
$ws.Range("Z9").Value = 'Old Mulford House'

# Row 10
$ws.Range("G10").Value = 'Rosenberg, Samuel'
$ws.Range("H10").Value = 626
$ws.Range("I10").Value = 'American'
$ws.Range("J10").Value = 1896
$ws.Range("K10").Value = 'Philadelphia (Philadelphia county, Pennsylvania, United States)'
$ws.Range("L10").Value = 'Abstract'
$ws.Range("N10").Value = 'Leopold, Otto Gerhard'
$ws.Range("O10").Value = 203
$ws.Range("P10").Value = 'American'
$ws.Range("Q10").Value = 1824
$ws.Range("R10").Value = 'Germany (Europe)'
$ws.Range("U10").Value = 'Cook, Robert A.'
$ws.Range("V10").Value = 141
$ws.Range("W10").Value = 'American'
$ws.Range("X10").Value = 1872
$ws.Range("Z10").Value = 'School; [floor plan] (en suite with 1997.29.34.1-.5)'

# Row 11
$ws.Range("G11").Value = 'Biddle, George'
$ws.Range("H11").Value = 87
$ws.Range("I11").Value = 'American'
$ws.Range("J11").Value = 1885
$ws.Range("L11").Value = 'Adam and Eve'
$ws.Range("M11").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=100726&size=Medium'
$ws.Range("N11").Value = 'Saint, Lawrence B.'
$ws.Range("O11").Value = 81
$ws.Range("P11").Value = 'American'
$ws.Range("Q11").Value = 1885
$ws.Range("R11").Value = 'Sharpsburgh, Pennsylvania'
$ws.Range("S11").Value = 'Drapery from a Sleeve of the Virgin, Window at West End of Church of St. Vincent, Rouen'
$ws.Range("T11").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=106095&size=Medium'
$ws.Range("U11").Value = 'Johnston, Ralph W.'
$ws.Range("V11").Value = 57
$ws.Range("W11").Value = 'American'
$ws.Range("X11").Value = 1873
$ws.Range("Z11").Value = '(Carnegie Library of Pittsburgh: Home Library Girl''s Club, February 7, 1907)'

# Row 12
$ws.Range("G12").Value = 'Strauss, Zoe'
$ws.Range("H12").Value = 221
$ws.Range("I12").Value = 'American'
$ws.Range("J12").Value = 1970
$ws.Range("K12").Value = 'Philadelphia, Pennsylvania, United States of America'
$ws.Range("L12").Value = 'Homesteading'
$ws.Range("M12").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=130922&size=Medium'
$ws.Range("N12").Value = 'Mellan, Claude'
$ws.Range("O12").Value = 120
$ws.Range("P12").Value = 'French'
$ws.Range("Q12").Value = 1598
$ws.Range("R12").Value = 'Abbeville, France'
$ws.Range("S12").Value = 'Antique Statue: Young Man'
$ws.Range("T12").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=105517&size=Medium'
$ws.Range("U12").Value = 'Mauve, Anton'
$ws.Range("V12").Value = 47
$ws.Range("W12").Value = 'Dutch'
$ws.Range("X12").Value = 1838
$ws.Range("Y12").Value = 'Zaandam, North Holland, Netherlands'
$ws.Range("Z12").Value = 'Wood Choppers'

# Row 13
$ws.Range("G13").Value = 'Bendiner, Alfred'
$ws.Range("H13").Value = 77
$ws.Range("I13").Value = 'American'
$ws.Range("J13").Value = 1899
$ws.Range("K13").Value = 'Pittsburgh, Pennsylvania'
$ws.Range("L13").Value = 'Travel is so Broadening'
$ws.Range("M13").Value = 'http://www.cmoa.org/CollectionImage.aspx?irn=89492&size=Medium|http://www.cmoa.org/CollectionImage.aspx?irn=47103&size=Medium'
$ws.Range("N13").Value = 'Walfish, Herbert S.'
$ws.Range("O13").Value = 73
$ws.Range("P13").Value = 'American'
$ws.Range("Q13").Value = 1923
$ws.Range("S13").Value = 'School; [exterior perspective] (en suite with 95.127.16-.22)'
$ws.Range("U13").Value = 'Deane, Edward E.'
$ws.Range("V13").Value = 54
$ws.Range("W13").Value = 'English'
$ws.Range("X13").Value = 1851
$ws.Range("Z13").Value = 'Untitled (Man), from Sketchbook (en suite with 91.23.12.1-.15)'

# Rows 5 and 6 (horoscopes "Snake" and "Horse") carry an explicit black font colour
# on the whole G:Z block in the source workbook (new fontId=3 / cellXfs index 1).
$ws.Range("G5:Z6").Font.Color = 0

# Restore the author's final selection from the saved workbook.
$ws.Range("L18").Select()
